$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 17500
$ws.Range("J13").Value = 17500
$ws.Range("L13").Value = 17500
$ws.Range("N13").Value = -17838

# Row 40
$ws.Range("H40").Value = 2086.7896
$ws.Range("I40").Value = 1849.875
$ws.Range("J40").Value = 2259.0908
$ws.Range("K40").Value = 1849.875
$ws.Range("L40").Value = 2259.0908
$ws.Range("M40").Value = -1674.875
$ws.Range("N40").Value = -2609.0908

# Row 113
$ws.Range("H113").Value = 2928.4285
$ws.Range("I113").Value = 2583.1667
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2583.1667
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 670.8332999999998
$ws.Range("N113").Value = -11508

# Row 125
$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 18000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -22920

# Row 135
$ws.Range("H135").Value = 1205.909
$ws.Range("I135").Value = 1092
$ws.Range("K135").Value = 9828
$ws.Range("M135").Value = -7293

# Row 137
$ws.Range("H137").Value = 1558.8182
$ws.Range("I137").Value = 1321.1428
$ws.Range("K137").Value = 3963.4284
$ws.Range("M137").Value = -1413.4284

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5268.1953
$ws.Range("I32").Value = 3135.027
$ws.Range("K32").Value = 3135.027
$ws.Range("M32").Value = -2848.027

# Row 35
$ws.Range("H35").Value = 1462.8334
$ws.Range("I35").Value = 1462.8334
$ws.Range("K35").Value = 1462.8334
$ws.Range("M35").Value = -1056.8334

# Row 61
$ws.Range("H61").Value = 1488.7778
$ws.Range("I61").Value = 1488.7778
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1488.7778
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1276.7778
$ws.Range("N61").ClearContents()

# Row 110
$ws.Range("H110").Value = 8738.9
$ws.Range("I110").Value = 9733.166999999999
$ws.Range("J110").Value = 7247.5
$ws.Range("K110").Value = 9733.166999999999
$ws.Range("L110").Value = 7247.5
$ws.Range("M110").Value = -7688.166999999999
$ws.Range("N110").Value = -11337.5

# Row 136
$ws.Range("H136").Value = 1488.7778
$ws.Range("I136").Value = 1488.7778
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4466.3334
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1916.3334
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 5406.75
$ws.Range("I20").Value = 6006
$ws.Range("K20").Value = 6006
$ws.Range("M20").Value = -5759

# Row 126
$ws.Range("H126").Value = 99888
$ws.Range("J126").Value = 99888
$ws.Range("L126").Value = 99888
$ws.Range("N126").Value = -109768

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 98.15000000000001
$ws.Range("I7").Value = 69.17646999999999
$ws.Range("J7").Value = 262.33334
$ws.Range("K7").Value = 69.17646999999999
$ws.Range("L7").Value = 262.33334
$ws.Range("M7").Value = 43.82353000000001
$ws.Range("N7").Value = -488.33334

# Row 99
$ws.Range("H99").Value = 17228.385
$ws.Range("I99").Value = 12247.5
$ws.Range("J99").Value = 19442.111
$ws.Range("K99").Value = 12247.5
$ws.Range("L99").Value = 19442.111
$ws.Range("M99").Value = -10749.5
$ws.Range("N99").Value = -22438.111

# Row 126
$ws.Range("H126").Value = 17228.385
$ws.Range("I126").Value = 12247.5
$ws.Range("J126").Value = 19442.111
$ws.Range("K126").Value = 36742.5
$ws.Range("L126").Value = 58326.333
$ws.Range("M126").Value = -34272.5
$ws.Range("N126").Value = -63266.333

# Row 132
$ws.Range("H132").Value = 1576.4
$ws.Range("I132").Value = 1580.421
$ws.Range("K132").Value = 4741.263
$ws.Range("M132").Value = -2211.263

# Row 134
$ws.Range("H134").Value = 2828.8572
$ws.Range("J134").Value = 2975
$ws.Range("L134").Value = 8925
$ws.Range("N134").Value = -13995

$ws = $wb.Worksheets.Item("CUL")
# Row 50
$ws.Range("H50").Value = 1458.2858
$ws.Range("J50").Value = 1999.6
$ws.Range("L50").Value = 5998.799999999999
$ws.Range("N50").Value = -6960.799999999999

# Row 53
$ws.Range("H53").Value = 1458.2858
$ws.Range("J53").Value = 1999.6
$ws.Range("L53").Value = 5998.799999999999
$ws.Range("N53").Value = -6960.799999999999

# Row 131
$ws.Range("H131").Value = 1326.1428
$ws.Range("I131").Value = 930
$ws.Range("J131").Value = 1333.3455
$ws.Range("K131").Value = 2790
$ws.Range("L131").Value = 4000.0365
$ws.Range("M131").Value = 2250
$ws.Range("N131").Value = -14080.0365

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()

# Row 123
$ws.Range("H123").Value = 49000.3
$ws.Range("J123").Value = 49000.3
$ws.Range("L123").Value = 49000.3
$ws.Range("N123").Value = -53900.3

# Row 126
$ws.Range("H126").Value = 4033.3
$ws.Range("J126").Value = 5006.75
$ws.Range("L126").Value = 15020.25
$ws.Range("N126").Value = -19960.25

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 4450.533
$ws.Range("J46").Value = 4821.4165
$ws.Range("L46").Value = 4821.4165
$ws.Range("N46").Value = -5197.4165

# Row 55
$ws.Range("H55").Value = 771
$ws.Range("I55").Value = 609
$ws.Range("K55").Value = 609
$ws.Range("M55").Value = -436

# Row 93
$ws.Range("H93").Value = 2320.2
$ws.Range("I93").Value = 2267
$ws.Range("J93").Value = 2400
$ws.Range("K93").Value = 2267
$ws.Range("L93").Value = 2400
$ws.Range("M93").Value = -1019
$ws.Range("N93").Value = -4896

# Row 100
$ws.Range("H100").Value = 2999
$ws.Range("I100").Value = 2999
$ws.Range("K100").Value = 2999
$ws.Range("M100").Value = -2458

# Row 132
$ws.Range("H132").Value = 3333.5293
$ws.Range("I132").Value = 3324.8
$ws.Range("K132").Value = 9974.400000000001
$ws.Range("M132").Value = -7444.400000000001

# Row 136
$ws.Range("H136").Value = 3108.5833
$ws.Range("J136").Value = 2004
$ws.Range("L136").Value = 6012
$ws.Range("N136").Value = -11112

$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("M75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("M78").ClearContents()

# Row 113
$ws.Range("H113").Value = 1030.3334
$ws.Range("J113").Value = 1148.3334
$ws.Range("L113").Value = 3445.0002
$ws.Range("N113").Value = -7785.0002

# Row 122
$ws.Range("H122").Value = 5801
$ws.Range("I122").Value = 7333.3335
$ws.Range("K122").Value = 22000.0005
$ws.Range("M122").Value = -19550.0005

# Row 126
$ws.Range("H126").Value = 6571.7144
$ws.Range("I126").Value = 5200.6
$ws.Range("K126").Value = 15601.8
$ws.Range("M126").Value = -13131.8

# Row 136
$ws.Range("H136").Value = 1470.2778
$ws.Range("I136").Value = 1470.2778
$ws.Range("K136").Value = 4410.8334
$ws.Range("M136").Value = -1860.8334

Write-Output "Applied all changes"
